$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.064.26'
$ws.Range("E2").Value = '  +4.97%  '

$ws.Range("D3").Value = '2.233.50'
$ws.Range("E3").Value = '  +4.05%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = "'261.12"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +3.42%  '

$ws.Range("D6").Value = "'83.34"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +14.89%  '

$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +3.30%  '

$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = "'0.604"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +5.13%  '

$ws.Range("D10").Value = "'44.47"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +12.78%  '

$ws.Range("D11").Value = "'0.0929"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +2.87%  '

$ws.Range("E12").Value = '  +5.65%  '

$ws.Range("E13").Value = '  +2.85%  '

$ws.Range("D14").Value = '2.567.29'
$ws.Range("E14").Value = '  +3.51%  '

$ws.Range("D15").Value = "'14.64"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +3.88%  '

$ws.Range("D16").Value = '2.233.25'
$ws.Range("E16").Value = '  +5.28%  '

$ws.Range("D17").Value = "'0.788"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +3.31%  '

$ws.Range("D18").Value = '43.961.46'
$ws.Range("E18").Value = '  +4.93%  '

$ws.Range("E19").Value = '  +3.07%  '

$ws.Range("D20").Value = "'71.14"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("D21").Value = "'6.04"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +4.43%  '

$ws.Range("D22").Value = "'2.39"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +12.41%  '

$ws.Range("D23").Value = "'232.52"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +3.31%  '

$ws.Range("D24").Value = "'9.26"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -2.54%  '

$ws.Range("D26").Value = "'10.82"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +4.16%  '

$ws.Range("D27").Value = "'40.67"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +11.49%  '

$ws.Range("D28").Value = "'3.36"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +1.91%  '

$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +3.04%  '

$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("D31").Value = "'172.94"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +2.95%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'20.71"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +4.44%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.0890"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +11.95%  '

$ws.Range("D34").Value = "'5.36"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +5.54%  '

$ws.Range("E35").Value = '  +9.21%  '

$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = "'0.123"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +2.59%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.0371"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +13.39%  '

$ws.Range("D38").Value = "'4.55"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +7.71%  '

$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").Value = "'13.15"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +10.62%  '

$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = "'3.03"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +28.17%  '

$ws.Range("E41").Value = '  +4.39%  '

$ws.Range("D42").Value = "'63.62"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +9.31%  '

$ws.Range("E43").Value = '  +8.97%  '

$ws.Range("E44").Value = '  +4.41%  '

$ws.Range("D45").Value = "'104.16"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +4.60%  '

$ws.Range("D46").Value = "'8.43"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +2.73%  '

$ws.Range("E47").Value = '  +3.02%  '

$ws.Range("E48").Value = '  +30.27%  '

$ws.Range("E49").Value = '  +5.37%  '

$ws.Range("D50").Value = "'0.445"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -3.70%  '

$ws.Range("E51").Value = '  +4.21%  '
